$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/value updates (Price column D) ---
$ws.Range("D2").Value  = "244.75"
$ws.Range("D4").Value  = "5.402"
$ws.Range("D5").Value  = "0.06037"
$ws.Range("D6").Value  = "3.393"
$ws.Range("D7").Value  = "0.8139"
$ws.Range("D8").Value  = "0.9226"
$ws.Range("D9").Value  = "0.1438"
$ws.Range("D10").Value = "0.07470"
$ws.Range("D11").Value = "0.03401"
$ws.Range("D12").Value = "0.03045"
$ws.Range("D13").Value = "0.09430"
$ws.Range("D16").Value = "0.04813"

# --- Rows 17-23: coins shifted by one position (ranking reshuffle) ---
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005532"
$ws.Range("E17").Value = "16TigerCashTCH"

$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D18").Value = "0.004165"
$ws.Range("E18").Value = "17HotbitTokenHTB"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Value = "0.0009893"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "LEO"
$ws.Range("C20").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D20").Value = "3.665"
$ws.Range("E20").Value = "19LEOLEO"

$ws.Range("B21").Value = "KuCoinToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D21").Value = "6.436"
$ws.Range("E21").Value = "20KuCoinTokenKCS"

$ws.Range("B22").Value = "BTSEToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D22").Value = "2.183"
$ws.Range("E22").Value = "21BTSETokenBTSE"

$ws.Range("B23").Value = "One"
$ws.Range("C23").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D23").Value = "0.01124"
$ws.Range("E23").Value = "22OneONEBestin24h"

# --- More simple price/value updates ---
$ws.Range("D26").Value = "0.00008401"
$ws.Range("D27").Value = "0.0002900"
$ws.Range("D40").Value = "0.04000"
$ws.Range("D41").Value = "0.006420"
$ws.Range("D42").Value = "0.1076"
$ws.Range("D43").Value = "0.002900"
$ws.Range("D44").Value = "0.005780"
$ws.Range("D45").Value = "0.00005254"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").Value = "0.002320"
